$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.125.96'
$ws.Range('E2').Value = '  +0.65%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.681.18'
$ws.Range('E3').Value = '  +0.34%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.20'
$ws.Range('E5').Value = '  +0.01%  '

# Row 6
$ws.Range('E6').Value = '  +0.36%  '

# Row 7
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('E8').Value = '  +2.10%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0623'
$ws.Range('E9').Value = '  +0.58%  '

# Row 10
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.29'
$ws.Range('E10').Value = '  +5.36%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.918.05'
$ws.Range('E12').Value = '  +0.37%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.670.53'
$ws.Range('E13').Value = '  +0.07%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.16'
$ws.Range('E14').Value = '  +1.74%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.537'
$ws.Range('E15').Value = '  +2.16%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.26'
$ws.Range('E16').Value = '  +0.92%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.117.82'
$ws.Range('E17').Value = '  +0.62%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '238.30'
$ws.Range('E18').Value = '  +1.33%  '

# Row 19
$ws.Range('E19').Value = '  +0.60%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0747'
$ws.Range('E20').Value = '  +1.88%  '

# Row 21
$ws.Range('E21').Value = '  +0.03%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.52'
$ws.Range('E22').Value = '  +1.64%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.46'
$ws.Range('E23').Value = '  +3.09%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.12'
$ws.Range('E24').Value = '  -2.89%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.75'
$ws.Range('E25').Value = '  +0.80%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.23'
$ws.Range('E26').Value = '  +0.81%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.33'
$ws.Range('E27').Value = '  +1.86%  '

# Row 28
$ws.Range('E28').Value = '  +0.75%  '

# Row 29
$ws.Range('E29').Value = '  +0.02%  '

# Row 30
$ws.Range('E30').Value = '  +0.32%  '

# Row 31
$ws.Range('E31').Value = '  +0.20%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.557.92'
$ws.Range('E32').Value = '  +5.34%  '

# Row 33
$ws.Range('E33').Value = '  +0.99%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.20'
$ws.Range('E34').Value = '  +1.69%  '

# Row 35
$ws.Range('E35').Value = '  +1.86%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.605'
$ws.Range('E36').Value = '  +3.76%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.936'
$ws.Range('E37').Value = '  +4.35%  '

# Row 38
$ws.Range('E38').Value = '  -1.27%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0175'
$ws.Range('E39').Value = '  +2.66%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.05'
$ws.Range('E40').Value = '  +0.86%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '69.20'
$ws.Range('E41').Value = '  +3.41%  '

# Row 42
$ws.Range('E42').Value = '  +0.04%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.67'
$ws.Range('E43').Value = '  -2.84%  '

# Row 44
$ws.Range('E44').Value = '  -1.71%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.825.65'
$ws.Range('E45').Value = '  +0.71%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.785'
$ws.Range('E46').Value = '  +0.46%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.87'
$ws.Range('E47').Value = '  +0.37%  '

# Row 48
$ws.Range('E48').Value = '  +3.20%  '

# Row 49
$ws.Range('E49').Value = '  +1.77%  '

# Row 50
$ws.Range('E50').Value = '  +3.29%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.10'
$ws.Range('E51').Value = '  +4.67%  '
